$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.910.94"
$ws.Cells.Item(2, 5).Value = "  +2.66%  "
$ws.Cells.Item(3, 4).Value = "3.067.61"
$ws.Cells.Item(3, 5).Value = "  +1.94%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "517.09"
$ws.Cells.Item(5, 5).Value = "  +1.50%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "142.99"
$ws.Cells.Item(6, 5).Value = "  +2.19%  "
$ws.Cells.Item(7, 5).Value = "  +0.06%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.437"
$ws.Cells.Item(8, 5).Value = "  +1.76%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "7.29"
$ws.Cells.Item(9, 5).Value = "  +2.56%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.108"
$ws.Cells.Item(10, 5).Value = "  -0.12%  "
$ws.Cells.Item(11, 5).Value = "  +2.22%  "
$ws.Cells.Item(12, 4).Value = "3.593.41"
$ws.Cells.Item(12, 5).Value = "  +2.31%  "
$ws.Cells.Item(13, 5).Value = "  +3.03%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "26.26"
$ws.Cells.Item(14, 5).Value = "  +3.40%  "
$ws.Cells.Item(15, 5).Value = "  +0.82%  "
$ws.Cells.Item(16, 4).Value = "57.922.19"
$ws.Cells.Item(16, 5).Value = "  +2.80%  "
$ws.Cells.Item(17, 4).Value = "3.059.33"
$ws.Cells.Item(17, 5).Value = "  +2.16%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.10"
$ws.Cells.Item(18, 5).Value = "  +2.75%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.84"
$ws.Cells.Item(19, 5).Value = "  -0.62%  "
$ws.Cells.Item(20, 5).Value = "  +1.10%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "332.61"
$ws.Cells.Item(21, 5).Value = "  +0.14%  "
$ws.Cells.Item(22, 5).Value = "  -0.25%  "
$ws.Cells.Item(23, 5).Value = "  +0.00%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "65.59"
$ws.Cells.Item(24, 5).Value = "  +1.68%  "
$ws.Cells.Item(25, 5).Value = "  +3.04%  "
$ws.Cells.Item(26, 5).Value = "  +0.03%  "
$ws.Cells.Item(27, 4).Value = "0.0₃0905"
$ws.Cells.Item(27, 5).Value = "  -3.65%  "
$ws.Cells.Item(28, 5).Value = "  +1.97%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.24"
$ws.Cells.Item(29, 5).Value = "  +5.16%  "
$ws.Cells.Item(30, 5).Value = "  +1.99%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.20"
$ws.Cells.Item(31, 5).Value = "  +2.57%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "20.70"
$ws.Cells.Item(32, 5).Value = "  +1.57%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "154.93"
$ws.Cells.Item(33, 5).Value = "  +1.77%  "
$ws.Cells.Item(34, 5).Value = "  +2.22%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.03"
$ws.Cells.Item(35, 5).Value = "  +3.44%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "27.02"
$ws.Cells.Item(36, 5).Value = "  +0.90%  "
$ws.Cells.Item(37, 5).Value = "  +4.57%  "
$ws.Cells.Item(38, 5).Value = "  +2.18%  "
$ws.Cells.Item(39, 4).Value = "3.108.53"
$ws.Cells.Item(39, 5).Value = "  +2.41%  "
$ws.Cells.Item(40, 5).Value = "  +3.74%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "36.46"
$ws.Cells.Item(41, 5).Value = "  +0.09%  "
$ws.Cells.Item(42, 5).Value = "  +0.05%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.655"
$ws.Cells.Item(43, 5).Value = "  +0.11%  "
$ws.Cells.Item(44, 4).Value = "2.263.93"
$ws.Cells.Item(44, 5).Value = "  +2.94%  "
$ws.Cells.Item(45, 5).Value = "  +8.24%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "20.79"
$ws.Cells.Item(46, 5).Value = "  +5.95%  "
$ws.Cells.Item(47, 5).Value = "  +1.93%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.945"
$ws.Cells.Item(48, 5).Value = "  +2.78%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "5.93"
$ws.Cells.Item(49, 5).Value = "  +1.47%  "
$ws.Cells.Item(50, 5).Value = "  +10.08%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "257.32"
$ws.Cells.Item(51, 5).Value = "  +12.21%  "

Write-Output "Updated cryptos list"